$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  35"
$ws.Range("C9").Value = "Report Covering the Week  8/28/2023  Through  9/3/2023"

# --- Precinct crime-stat table updates (rows 15-29) ---
$ws.Range("G15").Value = 1
$ws.Range("M15").Value = -37.5
$ws.Range("N15").Value = -76.744186046511
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 25
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 127
$ws.Range("J16").Value = 140
$ws.Range("K16").Value = -9.285714285714
$ws.Range("L16").Value = 4.95867768595
$ws.Range("M16").Value = -20.125786163522
$ws.Range("N16").Value = -79.180327868852
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -12.5
$ws.Range("F17").Value = 30
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 266
$ws.Range("J17").Value = 257
$ws.Range("K17").Value = 3.501945525291
$ws.Range("L17").Value = 9.465020576131
$ws.Range("M17").Value = 83.448275862069
$ws.Range("N17").Value = -36.515513126491
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -44.444444444444
$ws.Range("I18").Value = 77
$ws.Range("J18").Value = 123
$ws.Range("K18").Value = -37.398373983739
$ws.Range("L18").Value = 10
$ws.Range("M18").Value = 24.193548387096
$ws.Range("N18").Value = -74.836601307189
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -11.111111111111
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = 45.16129032258
$ws.Range("I19").Value = 283
$ws.Range("J19").Value = 315
$ws.Range("K19").Value = -10.15873015873
$ws.Range("L19").Value = 18.90756302521
$ws.Range("M19").Value = 55.494505494505
$ws.Range("N19").Value = -29.25
$ws.Range("C20").Value = 2
$ws.Range("D20").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -12.5
$ws.Range("I20").Value = 56
$ws.Range("J20").Value = 60
$ws.Range("K20").Value = -6.666666666666
$ws.Range("L20").Value = 36.585365853658
$ws.Range("M20").Value = 93.103448275862
$ws.Range("N20").Value = -79.56204379562
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 104
$ws.Range("G21").Value = 104
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 825
$ws.Range("J21").Value = 913
$ws.Range("K21").Value = -9.638554216867
$ws.Range("L21").Value = 12.551159618008
$ws.Range("M21").Value = 38.422818791946
$ws.Range("N21").Value = -60.202604920405
$ws.Range("D22").Value = 1
$ws.Range("F22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("E22").Value = -100
$ws.Range("H23").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("G22").Value = 1
$ws.Range("F22").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("H22").Value = 0
$ws.Range("H23").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("J22").Value = 4
$ws.Range("K22").Value = -25
$ws.Range("L22").Value = -25
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 8
$ws.Range("E23").Value = 12.5
$ws.Range("F23").Value = 32
$ws.Range("G23").Value = 41
$ws.Range("H23").Value = -21.951219512195
$ws.Range("I23").Value = 257
$ws.Range("J23").Value = 284
$ws.Range("K23").Value = -9.507042253521
$ws.Range("L23").Value = -7.885304659498
$ws.Range("M23").Value = 46.857142857142
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = -46.153846153846
$ws.Range("F24").Value = 67
$ws.Range("G24").Value = 89
$ws.Range("H24").Value = -24.719101123595
$ws.Range("I24").Value = 616
$ws.Range("J24").Value = 584
$ws.Range("K24").Value = 5.479452054794
$ws.Range("L24").Value = 21.021611001964
$ws.Range("M24").Value = 42.592592592592
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 45.454545454545
$ws.Range("F25").Value = 53
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = 29.268292682926
$ws.Range("I25").Value = 413
$ws.Range("J25").Value = 379
$ws.Range("K25").Value = 8.970976253298
$ws.Range("L25").Value = 19.020172910662
$ws.Range("M25").Value = -12.5
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C26").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = -57.142857142857
$ws.Range("J27").Value = 44
$ws.Range("K27").Value = -13.636363636363
$ws.Range("L27").Value = -25.490196078431
$ws.Range("C28").Value = 1
$ws.Range("G28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("F28").Value = 1
$ws.Range("G28").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("H28").Value = -66.666666666666
$ws.Range("I28").Value = 13
$ws.Range("K28").Value = -35
$ws.Range("L28").Value = -48
$ws.Range("M28").Value = -53.571428571428
$ws.Range("N28").Value = -81.159420289855
$ws.Range("C29").Value = 1
$ws.Range("G29").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("F29").Value = 1
$ws.Range("G29").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("H29").Value = -50
$ws.Range("I29").Value = 12
$ws.Range("K29").Value = -7.692307692307
$ws.Range("L29").Value = -47.826086956521
$ws.Range("M29").Value = -52
$ws.Range("N29").Value = -80.95238095238
